# Applies the weekly rotation of Fruta/Hortaliza data (Femacal de La Calera - Breva).
# The underlying data for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado) and S (Precio $/Kg) is
# shuffled across rows 2-9 (other columns stay identical).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D, M, N, O, P, S
$data = @{
    2 = @{ D = 44188; M = 30; N = 15000; O = 15000; P = 15000; S = 3000 }
    3 = @{ D = 44181; M = 30; N = 20000; O = 20000; P = 20000; S = 4000 }
    4 = @{ D = 44196; M = 56; N = 15000; O = 15000; P = 15000; S = 3000 }
    5 = @{ D = 44179; M = 45; N = 20000; O = 20000; P = 20000; S = 4000 }
    6 = @{ D = 44193; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    7 = @{ D = 44175; M = 25; N = 20000; O = 20000; P = 20000; S = 4000 }
    8 = @{ D = 44186; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
    9 = @{ D = 44189; M = 40; N = 15000; O = 15000; P = 15000; S = 3000 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
